$d = $word.ActiveDocument

# Locate the paragraph ending in Medhat's bio text, then remove every
# (empty) paragraph that follows it, leaving that paragraph as the
# document's last paragraph (immediately before the section break).
$anchorText = "Medhat is an instructor at BCIT, founder of the .NET BC Meetup Group in Vancouver, and a Microsoft MVP since 2017."

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $anchorText) {
        $anchorIndex = $i
    }
}

if ($anchorIndex -ge 1 -and $anchorIndex -lt $d.Paragraphs.Count) {
    $anchorEnd = $d.Paragraphs.Item($anchorIndex).Range.End
    $docEnd = $d.Content.End
    $trailing = $d.Range($anchorEnd, $docEnd)
    $trailing.Delete()
}
